$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FileFastq")

# This new pairing example illustrates a pairing type that does not work
# when only one side of the pair is present.

# Row 14 (f6 example): clear the relationship type, since with this new
# example it is now shown as a case that fails when one-sided.
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = $null

# Row 16: first side of the new pairing example (f7), intentionally left
# without a relationship type / related file to demonstrate the failure.
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "test_lab:f7_1, test_lab:alt_f7_1"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "fastq"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1"

# Row 17: second side of the new pairing example (f7), with relationship
# type and related file filled in.
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "test_lab:f7_2, test_lab:alt_f7_2"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "fastq"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "paired with"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "test_lab:alt_f7_1"

# Update the active selection to match the edited workbook
$ws.Range("F17").Select()

$wb.Save()
